$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "how are you"
$ws.Range("D1").Select()
